# Commit: "made chagnes for month append"
# A second ("May-24") reporting block (Cost %, Qty, Value, Rate) is appended
# in columns F:I, mirroring the existing "Apr-24" block in columns B:E, and
# two Apr-24 figures (HD Sale / Monofil Sales rows) are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths: widen col A (Particulars) and give the 4 new columns the
# same width already used by columns C-F.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20
$ws.Columns.Item(7).ColumnWidth = 15
$ws.Columns.Item(8).ColumnWidth = 15
$ws.Columns.Item(9).ColumnWidth = 15
$ws.Columns.Item(10).ColumnWidth = 15

# ---------------------------------------------------------------------------
# Columns F and I carry text-like values (percentages / rates-as-text), like
# columns B and E. Force text format first so values are stored as literal
# strings instead of being auto-parsed into numbers/percentages.
# ---------------------------------------------------------------------------
$ws.Range("F1:F35").NumberFormat = "@"
$ws.Range("I1:I35").NumberFormat = "@"

# ---------------------------------------------------------------------------
# New "May-24" block: F=Cost %, G=Qty, H=Value, I=Rate (rows 1-35).
# A lone "'" marks an intentionally-empty text cell (matches every blank cell
# already in this sheet, which is stored as an empty string, not left absent).
# ---------------------------------------------------------------------------
$colF = @("May-24", "Cost %", "55%", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "0%", "0%", "0%", "2%", "13%", "15%", "0%", "0%", "0%", "0%", "'", "0%", "0%", "0%", "0%", "0%", "0%", "0%", "15%", "'")
$colG = @("'", "Qty", 176401, 142379.9, 7416, 24415, 193921, 54225, 7896, 24415, 81691, "'", "'", 0, 4134, "'", "'", "'", 8271.5, 52769.8, "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'")
$colH = @("'", "Value", 17981704, 13105538.64, 1010633, 4649643, 19471496, 5235425, 1242959, 6078358, 20391750, 32948492, "'", 0, 41340, 11, "'", "'", 616661, 4168133, "'", "'", "'", "'", "'", "'", 125000, "'", "'", "'", "'", "'", "'", "'", "'")
$colI = @("'", "Rate", "101.94", "92.05", "136.28", "190.44", "100.41", "96.55", "157.42", "248.96", "249.62", "'", "'", "'", "10.00", "'", "'", "'", "74.55", "78.99", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'", "'")

for ($i = 0; $i -lt 35; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 6).Value = $colF[$i]
    $ws.Cells.Item($r, 7).Value = $colG[$i]
    $ws.Cells.Item($r, 8).Value = $colH[$i]
    $ws.Cells.Item($r, 9).Value = $colI[$i]
}

# ---------------------------------------------------------------------------
# Corrections to existing Apr-24 figures (rows 8 and 11). Columns E8/E11 are
# text-as-number "Rate" cells, so force text format before assigning.
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = 37950
$ws.Range("D8").Value = 3698350
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "97.45"
$ws.Range("C11").Value = 60566
$ws.Range("D11").Value = 14830233
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "244.86"

